$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K data to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats/styles from column F (old D, now shifted) into new D:E columns,
# restricted to the row blocks that actually contain quarterly data.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 4137900
$ws.Range("E8").Value = 4291900
$ws.Range("D9").Value = 3424100
$ws.Range("E9").Value = 3597900
$ws.Range("D10").Value = 713800
$ws.Range("E10").Value = 694000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 3882400
$ws.Range("E17").Value = 4045900
$ws.Range("D18").Value = 255500
$ws.Range("E18").Value = 246000
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 279800
$ws.Range("E21").Value = 269900
$ws.Range("D22").Value = 9500
$ws.Range("E22").Value = 6500
$ws.Range("D23").Value = 246100
$ws.Range("E23").Value = 239400
$ws.Range("D24").Value = 59100
$ws.Range("E24").Value = 61000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 187000
$ws.Range("E26").Value = 178500
$ws.Range("D27").Value = 187000
$ws.Range("E27").Value = 178500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 200
$ws.Range("E29").Value = -2600
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 187200
$ws.Range("E33").Value = 175900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 187200
$ws.Range("E35").Value = 175900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 378600
$ws.Range("E41").Value = 297800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 2322100
$ws.Range("E43").Value = 2453400
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 52400
$ws.Range("E45").Value = 53900
$ws.Range("D46").Value = 2753100
$ws.Range("E46").Value = 2805100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 228300
$ws.Range("E48").Value = 232000
$ws.Range("D49").Value = 1367700
$ws.Range("E49").Value = 1385400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 78300
$ws.Range("E52").Value = 70300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4427400
$ws.Range("E54").Value = 4492700
$ws.Range("D57").Value = 971000
$ws.Range("E57").Value = 1077800
$ws.Range("D58").Value = 97100
$ws.Range("E58").Value = 74900
$ws.Range("D59").Value = 365200
$ws.Range("E59").Value = 362700
$ws.Range("D60").Value = 1433300
$ws.Range("E60").Value = 1515300
$ws.Range("D61").Value = 1341400
$ws.Range("E61").Value = 1341300
$ws.Range("D62").Value = 57700
$ws.Range("E62").Value = 69500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2832300
$ws.Range("E66").Value = 2926100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 3845600
$ws.Range("E72").Value = 3728500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1595100
$ws.Range("E76").Value = 1566600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 187200
$ws.Range("E81").Value = 175900
$ws.Range("D83").Value = 24300
$ws.Range("E83").Value = 23900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 264000
$ws.Range("E89").Value = 220400
$ws.Range("D91").Value = -9200
$ws.Range("E91").Value = -15200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -20300
$ws.Range("E94").Value = -19600
$ws.Range("D96").Value = -70100
$ws.Range("E96").Value = -64600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -159800
$ws.Range("E100").Value = -204300
$ws.Range("D101").Value = -3100
$ws.Range("E101").Value = -9300
$ws.Range("D102").Value = 80800
$ws.Range("E102").Value = -12800
$ws.Range("F91").Value = -8900
$ws.Range("G91").Value = -11700
$ws.Range("H91").Value = -8000
$ws.Range("I91").Value = -8000
$ws.Range("J91").Value = -10600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F91").Value = -8900
$ws.Range("G91").Value = -11700
$ws.Range("H91").Value = -8000
$ws.Range("I91").Value = -8000
$ws.Range("J91").Value = -10600

